# Insert a new price-record row at row 40 of the weekly Membrillo sheet.
# This pushes the existing rows 40..60 down to 41..61 (Excel copies the
# formatting of the row above when inserting, which matches the target
# workbook's formatting for the date column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 40.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new weekly record.
$ws.Cells.Item(40, 1).Value  = 9
$ws.Cells.Item(40, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(40, 3).Value  = "Metropolitana"
$ws.Cells.Item(40, 4).Value  = 45029
$ws.Cells.Item(40, 5).Value  = 13
$ws.Cells.Item(40, 6).Value  = "Fruta"
$ws.Cells.Item(40, 7).Value  = 100104
$ws.Cells.Item(40, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(40, 9).Value  = 100104003
$ws.Cells.Item(40, 10).Value = "Membrillo"
$ws.Cells.Item(40, 11).Value = "Champion"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 18
$ws.Cells.Item(40, 14).Value = 225000
$ws.Cells.Item(40, 15).Value = 270000
$ws.Cells.Item(40, 16).Value = 250000
$ws.Cells.Item(40, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(40, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(40, 19).Value = 556
$ws.Cells.Item(40, 20).Value = 450
